# Add a new "highest_reading_note" column (K) to the musical_instruments
# table, mirroring the existing "lowest_reading_note" column (J):
# most instruments get "NA", but double_bass and electric_bass get 64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the new column similar to the other "reading note" column.
$ws.Columns.Item(11).ColumnWidth = 38.83

# Header cell - bold, matching the other header cells in row 1.
$ws.Range("K1").Font.Name = "Arial"
$ws.Range("K1").Font.Size = 12
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Value = "highest_reading_note"

$values = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", 64, "NA", 64, "NA")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 11)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 12
    $cell.Font.Color = 0
    $cell.Value = $values[$i]
}

# Match the new selection noted in the saved workbook view.
$ws.Range("K18").Select()
